$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of chat data below the existing two data rows.
$ws.Range("A3").Value = "14-12-2021 12:00"
$ws.Range("B3").Value = ""

$ws.Range("A4").Value = "14-12-2021 11:16"
$ws.Range("B4").Value = "Para los que no sabían, tengo dos hijos y bien chulos🥴💗"

$ws.Range("A5").Value = "14-12-2021 11:28"
$ws.Range("B5").Value = "Mi pequeño unicornio 🥰💗"

# Match the pageSetup attributes added in the target workbook.
$ws.PageSetup.FirstPageNumber = 1
$ws.PageSetup.UseFirstPageNumber = $true
$ws.PageSetup.Copies = 1
